$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.9
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 2.4
$ws.Range("L2").Value = 3.75
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 21
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 9
$ws.Range("AC2").Value = 21
$ws.Range("AD2").Value = 8.5
$ws.Range("AL2").Value = 23
$ws.Range("AO2").Value = 9.5
$ws.Range("AQ2").Value = 29
$ws.Range("BB2").Value = 51
